# Auto-generated script applying numeric corrections to H:N leve-profit columns
# across multiple worksheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2279.5881
$ws.Range("I28").Value = 1800.4667
$ws.Range("K28").Value = 1800.4667
$ws.Range("M28").Value = -1315.4667
$ws.Range("H38").Value = 1212.0625
$ws.Range("J38").Value = 2485.2856
$ws.Range("L38").Value = 7455.8568
$ws.Range("N38").Value = -8199.856800000001
$ws.Range("H62").Value = 22731790
$ws.Range("I62").Value = 35718668
$ws.Range("K62").Value = 35718668
$ws.Range("M62").Value = -35718044
$ws.Range("H65").Value = 22731790
$ws.Range("I65").Value = 35718668
$ws.Range("K65").Value = 178593340
$ws.Range("M65").Value = -178590220
$ws.Range("H74").Value = 5808
$ws.Range("J74").Value = 6151.773
$ws.Range("L74").Value = 6151.773
$ws.Range("N74").Value = -8023.773
$ws.Range("H77").Value = 5808
$ws.Range("J77").Value = 6151.773
$ws.Range("L77").Value = 30758.865
$ws.Range("N77").Value = -40118.86500000001
$ws.Range("H100").Value = 8242.933999999999
$ws.Range("I100").Value = 2274.5
$ws.Range("K100").Value = 2274.5
$ws.Range("M100").Value = -1733.5
$ws.Range("H107").Value = 805.8570999999999
$ws.Range("I107").Value = 778.4
$ws.Range("K107").Value = 778.4
$ws.Range("M107").Value = 1141.6
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H131").Value = 5432
$ws.Range("I131").Value = 2726.8572
$ws.Range("K131").Value = 8180.571599999999
$ws.Range("M131").Value = -3140.571599999999
$ws.Range("H132").Value = 3306.7021
$ws.Range("I132").Value = 3232.111
$ws.Range("J132").Value = 4985
$ws.Range("K132").Value = 9696.332999999999
$ws.Range("L132").Value = 14955
$ws.Range("M132").Value = -7166.332999999999
$ws.Range("N132").Value = -20015
$ws.Range("H137").Value = 1708.4736
$ws.Range("I137").Value = 1301.9
$ws.Range("J137").Value = 2160.2222
$ws.Range("K137").Value = 3905.7
$ws.Range("L137").Value = 6480.6666
$ws.Range("M137").Value = -1355.7
$ws.Range("N137").Value = -11580.6666
$ws.Range("H138").Value = 3369.7727
$ws.Range("J138").Value = 3134.2144
$ws.Range("L138").Value = 9402.643199999999
$ws.Range("N138").Value = -19682.6432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3214.1538
$ws.Range("I32").Value = 2315.6943
$ws.Range("K32").Value = 2315.6943
$ws.Range("M32").Value = -2028.6943
$ws.Range("H74").Value = 2283.7083
$ws.Range("I74").Value = 2163
$ws.Range("J74").Value = 3128.6667
$ws.Range("K74").Value = 2163
$ws.Range("L74").Value = 3128.6667
$ws.Range("M74").Value = -1289
$ws.Range("N74").Value = -4876.6667
$ws.Range("H77").Value = 2283.7083
$ws.Range("I77").Value = 2163
$ws.Range("J77").Value = 3128.6667
$ws.Range("K77").Value = 10815
$ws.Range("L77").Value = 15643.3335
$ws.Range("M77").Value = -6447
$ws.Range("N77").Value = -24379.3335
$ws.Range("H132").Value = 1979.6428
$ws.Range("I132").Value = 1841.4166
$ws.Range("J132").Value = 2809
$ws.Range("K132").Value = 5524.2498
$ws.Range("L132").Value = 8427
$ws.Range("M132").Value = -2994.2498
$ws.Range("N132").Value = -13487

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 101500
$ws.Range("J57").Value = 101500
$ws.Range("L57").Value = 101500
$ws.Range("N57").Value = -102940
$ws.Range("H94").Value = 88333.336
$ws.Range("I94").Value = 65000
$ws.Range("J94").Value = 100000
$ws.Range("K94").Value = 65000
$ws.Range("L94").Value = 100000
$ws.Range("M94").Value = -64549
$ws.Range("N94").Value = -100902
$ws.Range("H136").Value = 101500
$ws.Range("J136").Value = 101500
$ws.Range("L136").Value = 101500
$ws.Range("N136").Value = -111700

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2113.5173
$ws.Range("I31").Value = 1757.3846
$ws.Range("J31").Value = 2402.875
$ws.Range("K31").Value = 1757.3846
$ws.Range("L31").Value = 2402.875
$ws.Range("M31").Value = -1462.3846
$ws.Range("N31").Value = -2992.875
$ws.Range("H34").Value = 2113.5173
$ws.Range("I34").Value = 1757.3846
$ws.Range("J34").Value = 2402.875
$ws.Range("K34").Value = 1757.3846
$ws.Range("L34").Value = 2402.875
$ws.Range("M34").Value = -1555.3846
$ws.Range("N34").Value = -2806.875
$ws.Range("H99").Value = 2500.5625
$ws.Range("I99").Value = 2453.9333
$ws.Range("K99").Value = 2453.9333
$ws.Range("M99").Value = -955.9333000000001
$ws.Range("H122").Value = 2253.1667
$ws.Range("I122").Value = 1755.125
$ws.Range("K122").Value = 5265.375
$ws.Range("M122").Value = -2815.375
$ws.Range("H126").Value = 2500.5625
$ws.Range("I126").Value = 2453.9333
$ws.Range("K126").Value = 7361.7999
$ws.Range("M126").Value = -4891.7999
$ws.Range("H132").Value = 3506.7715
$ws.Range("I132").Value = 3370.4285
$ws.Range("J132").Value = 4052.1428
$ws.Range("K132").Value = 10111.2855
$ws.Range("L132").Value = 12156.4284
$ws.Range("M132").Value = -7581.2855
$ws.Range("N132").Value = -17216.4284
$ws.Range("H134").Value = 1386.1904
$ws.Range("I134").Value = 1205.5
$ws.Range("K134").Value = 3616.5
$ws.Range("M134").Value = -1081.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 37049944
$ws.Range("I138").Value = 100003730
$ws.Range("K138").Value = 300011190
$ws.Range("M138").Value = -300006050

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14251.833
$ws.Range("I70").Value = 17628.25
$ws.Range("J70").Value = 7499
$ws.Range("K70").Value = 17628.25
$ws.Range("L70").Value = 7499
$ws.Range("M70").Value = -17358.25
$ws.Range("N70").Value = -8039
$ws.Range("H73").Value = 14251.833
$ws.Range("I73").Value = 17628.25
$ws.Range("J73").Value = 7499
$ws.Range("K73").Value = 17628.25
$ws.Range("L73").Value = 7499
$ws.Range("M73").Value = -16692.25
$ws.Range("N73").Value = -9371
$ws.Range("H132").Value = 3007.3
$ws.Range("I132").Value = 3007.3
$ws.Range("K132").Value = 9021.900000000001
$ws.Range("M132").Value = -6491.900000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2673.2307
$ws.Range("I40").Value = 2556
$ws.Range("K40").Value = 2556
$ws.Range("M40").Value = -2420
$ws.Range("H46").Value = 1965.9166
$ws.Range("I46").Value = 1257
$ws.Range("J46").Value = 2958.4
$ws.Range("K46").Value = 1257
$ws.Range("L46").Value = 2958.4
$ws.Range("M46").Value = -1069
$ws.Range("N46").Value = -3334.4
$ws.Range("H55").Value = 787.3214
$ws.Range("I55").Value = 197.3
$ws.Range("K55").Value = 197.3
$ws.Range("M55").Value = -24.30000000000001
$ws.Range("H61").Value = 5126.75
$ws.Range("I61").Value = 5126.75
$ws.Range("K61").Value = 5126.75
$ws.Range("M61").Value = -4924.75
$ws.Range("H113").Value = 5126.75
$ws.Range("I113").Value = 5126.75
$ws.Range("K113").Value = 5126.75
$ws.Range("M113").Value = -2956.75
$ws.Range("H132").Value = 2707.5247
$ws.Range("I132").Value = 2306.111
$ws.Range("J132").Value = 3285.56
$ws.Range("K132").Value = 6918.333
$ws.Range("L132").Value = 9856.68
$ws.Range("M132").Value = -4388.333
$ws.Range("N132").Value = -14916.68
$ws.Range("H136").Value = 2473.175
$ws.Range("I136").Value = 2200.92
$ws.Range("J136").Value = 2926.9333
$ws.Range("K136").Value = 6602.76
$ws.Range("L136").Value = 8780.7999
$ws.Range("M136").Value = -4052.76
$ws.Range("N136").Value = -13880.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 11749.5
$ws.Range("I2").Value = 11749.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 11749.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -11637.5
$ws.Range("N2").ClearContents()
$ws.Range("H96").Value = 842.9091
$ws.Range("I96").Value = 835.7778
$ws.Range("K96").Value = 835.7778
$ws.Range("M96").Value = 537.2222
$ws.Range("H122").Value = 6098.241
$ws.Range("I122").Value = 6098.241
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 18294.723
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15844.723
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 3301.2222
$ws.Range("I132").Value = 3030.182
$ws.Range("J132").Value = 3727.1428
$ws.Range("K132").Value = 9090.545999999998
$ws.Range("L132").Value = 11181.4284
$ws.Range("M132").Value = -6560.545999999998
$ws.Range("N132").Value = -16241.4284
$ws.Range("H136").Value = 4394.029
$ws.Range("I136").Value = 2659.7188
$ws.Range("K136").Value = 7979.1564
$ws.Range("M136").Value = -5429.1564
